# Generate Report for Handback
# Adds a new handback row (e2218a0a-fcf8-4b9d-99da-176901c1c0ad.md) to the
# Overview / zh-cn / de-de sheets, alongside the already-existing row for
# b21e5dab-a5f3-426c-8ff9-0d424568832d.md (renamed from the former
# 7a76190d-de93-40ea-8880-55568815d466.md).

$wb = $excel.ActiveWorkbook

$oldGuid = "7a76190d-de93-40ea-8880-55568815d466"
$renamedGuid = "b21e5dab-a5f3-426c-8ff9-0d424568832d"
$newGuid = "e2218a0a-fcf8-4b9d-99da-176901c1c0ad"

$renamedFile = "$renamedGuid.md"
$newFile = "$newGuid.md"

$zhXlfHash = "1f07cb62c21f1865c6b2d9c37dbc3adaab67ec21"
$newXlfHash = "df638fde329d85561d507f8989b40c71c17a3bb4"

$renamedZhXlf = "$renamedGuid.$zhXlfHash.zh-cn.xlf"
$renamedDeXlf = "$renamedGuid.$zhXlfHash.de-de.xlf"
$newZhXlf = "$newGuid.$newXlfHash.zh-cn.xlf"
$newDeXlf = "$newGuid.$newXlfHash.de-de.xlf"

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet "Overview" - rename the first data row's displayed hyperlink text
# (file name changed from the old guid to the "renamed" guid) and append
# a brand new row describing the second handback file.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)

# existing row 2 hyperlink text changes to reflect the renamed file
$wsOverview.Hyperlinks.Item(1).TextToDisplay = "e2e\$renamedFile"

$rowOv = $loOverview.ListRows.Add()
$wsOverview.Range("A3").Value = $newFile
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G3").Value = "2016-08-31 07:09:51"
$wsOverview.Range("G3").NumberFormat = $dateFmt
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6e137cb76b1c0ff7456170f8c7a5e471d38acfdb/e2e/$newFile", "", "", "e2e\$newFile") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)

$wsZh.Hyperlinks.Item(1).TextToDisplay = $renamedFile
$wsZh.Hyperlinks.Item(2).TextToDisplay = $renamedFile
$wsZh.Range("G2").Value = $renamedZhXlf
$wsZh.Range("H2").Value = "2016-08-31 07:09:46"
$wsZh.Range("J2").Value = $renamedZhXlf
$wsZh.Range("K2").Value = "2016-08-31 07:10:15"

$rowZh = $loZh.ListRows.Add()
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "True"
$wsZh.Range("G3").Value = $newZhXlf
$wsZh.Range("H3").Value = "2016-08-31 07:09:46"
$wsZh.Range("H3").NumberFormat = $dateFmt
$wsZh.Range("J3").Value = $newZhXlf
$wsZh.Range("K3").Value = "2016-08-31 07:10:15"
$wsZh.Range("K3").NumberFormat = $dateFmt
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = "True"
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "False"
$wsZh.Range("P3").Value = ""
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6e137cb76b1c0ff7456170f8c7a5e471d38acfdb/e2e/$newFile", "", "", $newFile) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/4f4dc182a5ecdc6cb2d92bd4500b274224c2e538/e2e/$newFile", "", "", $newFile) | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)

$wsDe.Hyperlinks.Item(1).TextToDisplay = $renamedFile
$wsDe.Hyperlinks.Item(2).TextToDisplay = $renamedFile
$wsDe.Range("G2").Value = $renamedDeXlf
$wsDe.Range("J2").Value = $renamedDeXlf
$wsDe.Range("K2").Value = "2016-08-31 07:10:28"

$rowDe = $loDe.ListRows.Add()
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "True"
$wsDe.Range("G3").Value = $newDeXlf
$wsDe.Range("H3").Value = "2016-08-31 07:09:51"
$wsDe.Range("H3").NumberFormat = $dateFmt
$wsDe.Range("J3").Value = $newDeXlf
$wsDe.Range("K3").Value = "2016-08-31 07:10:28"
$wsDe.Range("K3").NumberFormat = $dateFmt
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = "True"
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "False"
$wsDe.Range("P3").Value = ""
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6e137cb76b1c0ff7456170f8c7a5e471d38acfdb/e2e/$newFile", "", "", $newFile) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/e23dd405ae28f7f3fa43c72685f3cab4ccaf4d6e/e2e/$newFile", "", "", $newFile) | Out-Null

Write-Output "Handback report row added."
